$d = $word.ActiveDocument

$replacements = @(
    @("2025-07-18 Friday", "2025-07-19 Saturday"),
    @("837÷9=", "608÷2="),
    @("768÷5=", "849÷6="),
    @("140÷3=", "116÷9="),
    @("878÷7=", "410÷7="),
    @("292÷6=", "489÷5="),
    @("634÷2=", "545÷9="),
    @("172÷2=", "459÷6="),
    @("767÷5=", "867÷5="),
    @("347÷7=", "637÷5="),
    @("833÷4=", "519÷3="),
    @("383÷8=", "318÷4="),
    @("326÷8=", "551÷3="),
    @("519÷9=", "261÷8="),
    @("429÷7=", "788÷4="),
    @("971÷7=", "623÷8="),
    @("540÷2=", "690÷8="),
    @("480÷2=", "312÷2="),
    @("209÷2=", "618÷7="),
    @("987÷6=", "582÷3="),
    @("479÷6=", "283÷5="),
    @("438÷3=", "859÷2="),
    @("499÷9=", "130÷3="),
    @("322÷9=", "879÷7="),
    @("907÷6=", "913÷7="),
    @("254÷3=", "518÷8=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
